$wb = $excel.ActiveWorkbook

# Helper: write a text label into a cell even when the label reads like a
# number (e.g. "2050"), without Excel re-interpreting it as a numeric value
# (and without bloating styles.xml with a new quote-prefixed style). We do
# this by building the label as a text formula result in a scratch cell,
# copying it, and pasting values-only on top of the target cell - the
# target cell's existing style/formatting is left completely untouched.
function Set-TextLabel($ws, [string]$cellRef, [string]$text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163) # xlPasteValues
    $scratch.Clear()
}

# Sheets whose "2050" (or "2041-2050") column header was wrongly stored as
# the stray numeric value 695.7894663120405 instead of the intended text
# label, and which carry a trailing "Total" row to drop.
$sheetConfigs = @(
    @{ Name = "Potencia Acumulada - SIN (MW)"; E1Label = "2050"; TotalRow = 13 },
    @{ Name = "Geracao Periodo Medio (MWMed)"; E1Label = "2050"; TotalRow = 13 },
    @{ Name = "Atendimento a Ponta(MW)"; E1Label = "2050"; TotalRow = 13 },
    @{ Name = "Potencia Incremental - SIN(MW)"; E1Label = "2041-2050"; TotalRow = 13 }
)

foreach ($cfg in $sheetConfigs) {
    $ws = $wb.Worksheets.Item($cfg.Name)
    Set-TextLabel $ws "E1" $cfg.E1Label
    $ws.Rows.Item($cfg.TotalRow).Delete()
}

# "Custo Total (bilhões de R$)" sheet only needs its Total row removed.
$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$wsCusto.Rows.Item(4).Delete()
